$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.458702551532514
$ws.Range("C2").Value = 6.6092802915852
$ws.Range("E2").Value = 16.53154759281661
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 21.59341235664857
$ws.Range("H2").Value = 12.55754310208166
$ws.Range("K2").Value = 8.606658067835472
$ws.Range("N2").Value = 17.01278109724849
$ws.Range("O2").Value = 18.1438398480341
$ws.Range("B3").Value = 8.115647334979764
$ws.Range("C3").Value = 6.504870617224146
$ws.Range("E3").Value = 15.59355093206447
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 21.67150609310042
$ws.Range("H3").Value = 12.60603450561432
$ws.Range("K3").Value = 8.318092797627614
$ws.Range("N3").Value = 17.06341470440659
$ws.Range("O3").Value = 18.22440881970333
$ws.Range("B4").Value = 7.898445033565831
$ws.Range("C4").Value = 6.43946914583015
$ws.Range("E4").Value = 14.99236570412395
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 21.72802924302474
$ws.Range("H4").Value = 12.63783454778404
$ws.Range("K4").Value = 8.134088273594539
$ws.Range("N4").Value = 17.09612299366481
$ws.Range("O4").Value = 18.27807516096775
$ws.Range("B5").Value = 7.808414423413815
$ws.Range("C5").Value = 6.412512703949927
$ws.Range("E5").Value = 14.74129623204272
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 21.75320349798038
$ws.Range("H5").Value = 12.65130285277915
$ws.Range("K5").Value = 8.057456980420845
$ws.Range("N5").Value = 17.1098600662367
$ws.Range("O5").Value = 18.30099725310113
$ws.Range("B6").Value = 7.793377335302663
$ws.Range("C6").Value = 6.408018807022661
$ws.Range("E6").Value = 14.69924790250193
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 21.75751249130493
$ws.Range("H6").Value = 12.65357003054764
$ws.Range("K6").Value = 8.044635012601892
$ws.Range("N6").Value = 17.11216578088592
$ws.Range("O6").Value = 18.30486694713226
$ws.Range("B7").Value = 7.897236813418629
$ws.Range("C7").Value = 6.439106808516519
$ws.Range("E7").Value = 14.98900392459326
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 21.72836010531663
$ws.Range("H7").Value = 12.63801412293196
$ws.Range("K7").Value = 8.133061374395504
$ws.Range("N7").Value = 17.09630660253028
$ws.Range("O7").Value = 18.27838003783309
$ws.Range("B8").Value = 8.341854027693083
$ws.Range("C8").Value = 6.573560375829725
$ws.Range("E8").Value = 16.2135124984723
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 21.61855079467341
$ws.Range("H8").Value = 12.57384244048867
$ws.Range("K8").Value = 8.508615742638394
$ws.Range("N8").Value = 17.02990416701719
$ws.Range("O8").Value = 18.17074741004442
$ws.Range("B9").Value = 9.156630448544028
$ws.Range("C9").Value = 6.826096337704097
$ws.Range("E9").Value = 18.50122265353804
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 21.4718802598642
$ws.Range("H9").Value = 12.46407639511772
$ws.Range("K9").Value = 9.188300156732025
$ws.Range("N9").Value = 16.91248731250109
$ws.Range("O9").Value = 17.99309801832483
$ws.Range("B10").Value = 9.714686021141297
$ws.Range("C10").Value = 7.003748045112257
$ws.Range("E10").Value = 20.14047895013242
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 21.40679416864066
$ws.Range("H10").Value = 12.39322806181529
$ws.Range("K10").Value = 9.650002450673695
$ws.Range("N10").Value = 16.83395541908365
$ws.Range("O10").Value = 17.88311470592888
$ws.Range("B11").Value = 9.958793733714884
$ws.Range("C11").Value = 7.082636250082874
$ws.Range("E11").Value = 20.84377544918225
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 21.38660309402866
$ws.Range("H11").Value = 12.3631244681813
$ws.Range("K11").Value = 9.851342237259935
$ws.Range("N11").Value = 16.799894366046
$ws.Range("O11").Value = 17.83757480582566
$ws.Range("B12").Value = 10.04976304100815
$ws.Range("C12").Value = 7.112215033205353
$ws.Range("E12").Value = 21.10403266994742
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 21.38032210572398
$ws.Range("H12").Value = 12.35203077960139
$ws.Range("K12").Value = 9.926298081115123
$ws.Range("N12").Value = 16.78723449373529
$ws.Range("O12").Value = 17.82097895446393
$ws.Range("B13").Value = 10.03023746027107
$ws.Range("C13").Value = 7.105858077661827
$ws.Range("E13").Value = 21.04825088802091
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 21.38161396461485
$ws.Range("H13").Value = 12.35440639569364
$ws.Range("K13").Value = 9.910212790471322
$ws.Range("N13").Value = 16.78995043903885
$ws.Range("O13").Value = 17.82452424777498
$ws.Range("B14").Value = 9.966307619575282
$ws.Range("C14").Value = 7.085075712638438
$ws.Range("E14").Value = 20.86530835161275
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 21.38605894066687
$ws.Range("H14").Value = 12.36220565193699
$ws.Range("K14").Value = 9.857534893456286
$ws.Range("N14").Value = 16.79884805975729
$ws.Range("O14").Value = 17.83619642767401
$ws.Range("B15").Value = 9.926955638303303
$ws.Range("C15").Value = 7.072307063321194
$ws.Range("E15").Value = 20.7524617675766
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 21.38895967049233
$ws.Range("H15").Value = 12.36702276426414
$ws.Range("K15").Value = 9.825099506122539
$ws.Range("N15").Value = 16.80432911782637
$ws.Range("O15").Value = 17.84343060349211
$ws.Range("B16").Value = 9.698530830907712
$ws.Range("C16").Value = 6.998552238146677
$ws.Range("E16").Value = 20.09366626968082
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 21.40830414610026
$ws.Range("H16").Value = 12.39523819470094
$ws.Range("K16").Value = 9.636666227664509
$ws.Range("N16").Value = 16.83621478556349
$ws.Range("O16").Value = 17.88618150807148
$ws.Range("B17").Value = 9.555850919541548
$ws.Range("C17").Value = 6.952799930855394
$ws.Range("E17").Value = 19.67867771883401
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 21.42259132289174
$ws.Range("H17").Value = 12.41309203223089
$ws.Range("K17").Value = 9.518814363279912
$ws.Range("N17").Value = 16.85620102657693
$ws.Range("O17").Value = 17.9135605822309
$ws.Range("B18").Value = 9.472870361806093
$ws.Range("C18").Value = 6.926304073873046
$ws.Range("E18").Value = 19.43599357666292
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 21.43169472439938
$ws.Range("H18").Value = 12.42356115077021
$ws.Range("K18").Value = 9.450212966992478
$ws.Range("N18").Value = 16.86785321293917
$ws.Range("O18").Value = 17.9297309260269
$ws.Range("B19").Value = 9.444619626274775
$ws.Range("C19").Value = 6.91730260487213
$ws.Range("E19").Value = 19.35313797865465
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 21.43492876692924
$ws.Range("H19").Value = 12.42714016915723
$ws.Range("K19").Value = 9.426846764603905
$ws.Range("N19").Value = 16.87182537193694
$ws.Range("O19").Value = 17.93527841624046
$ws.Range("B20").Value = 9.571134633908732
$ws.Range("C20").Value = 6.957689135859781
$ws.Range("E20").Value = 19.72326705547632
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 21.42097866988018
$ws.Range("H20").Value = 12.41117075338955
$ws.Range("K20").Value = 9.531444647487517
$ws.Range("N20").Value = 16.85405725313009
$ws.Range("O20").Value = 17.91060227001558
$ws.Range("B21").Value = 9.98512571196248
$ws.Range("C21").Value = 7.09118812331372
$ws.Range("E21").Value = 20.91920737792269
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 21.38471621973124
$ws.Range("H21").Value = 12.35990651771918
$ws.Range("K21").Value = 9.873042876212892
$ws.Range("N21").Value = 16.79622815231721
$ws.Range("O21").Value = 17.83275038332618
$ws.Range("B22").Value = 10.2471062527657
$ws.Range("C22").Value = 7.176713598191527
$ws.Range("E22").Value = 21.6654817288934
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 21.36897660078095
$ws.Range("H22").Value = 12.32818544873265
$ws.Range("K22").Value = 10.0887804582766
$ws.Range("N22").Value = 16.75982213451941
$ws.Range("O22").Value = 17.78565471249262
$ws.Range("B23").Value = 10.10808744309288
$ws.Range("C23").Value = 7.131230300652769
$ws.Range("E23").Value = 21.27040415540426
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 21.37664555237615
$ws.Range("H23").Value = 12.34495235067443
$ws.Range("K23").Value = 9.974336152056711
$ws.Range("N23").Value = 16.77912594365121
$ws.Range("O23").Value = 17.81044318162983
$ws.Range("B24").Value = 9.56422782669158
$ws.Range("C24").Value = 6.955479323787466
$ws.Range("E24").Value = 19.70312099518062
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 21.42170498030524
$ws.Range("H24").Value = 12.4120387258864
$ws.Range("K24").Value = 9.525737130718566
$ws.Range("N24").Value = 16.85502594892536
$ws.Range("O24").Value = 17.91193838451956
$ws.Range("B25").Value = 8.94297051127996
$ws.Range("C25").Value = 6.75908682904968
$ws.Range("E25").Value = 17.85972769427484
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 21.50411958619387
$ws.Range("H25").Value = 12.49205049182497
$ws.Range("K25").Value = 9.010845015177578
$ws.Range("N25").Value = 16.9428886250038
$ws.Range("O25").Value = 18.03756244161513
